$d = $word.ActiveDocument

# Locate the existing " IDEA 2018.3.5, " run text (right after "IntelliJ").
$target = $d.Content
$found = $target.Find.Execute(" IDEA 2018.3.5, ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rangeStart = $target.Start

    # Split point: right after " IDEA 2018.3.5" (before the ", ").
    $splitAt = $rangeStart + (" IDEA 2018.3.5").Length

    # Insert the new build-version text right at the split point; it
    # inherits the surrounding (Bold, Navy, sz20) direct formatting.
    $insertionPoint = $d.Range($splitAt, $splitAt)
    $insertionPoint.InsertAfter(" IU 183.5912.21")

    # The newly inserted text ends up merged into the same run as its
    # neighbours because formatting is identical. Briefly bookmarking
    # exactly the inserted span forces Word to keep it as a distinct
    # run (bookmark boundaries cannot sit inside a run); removing the
    # bookmark afterwards leaves three separate, identically-formatted
    # runs behind, matching a manual "insert new run" edit.
    $newTextLen = (" IU 183.5912.21").Length
    $newRunRange = $d.Range($splitAt, $splitAt + $newTextLen)
    $d.Bookmarks.Add("TempBuildVersionSplit", $newRunRange) | Out-Null
    $d.Bookmarks("TempBuildVersionSplit").Delete()
}
